$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.404.65"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.846.74"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "'0.9982"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'240.17"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'0.6361"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07558"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.2969"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'24.68"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").Value = "'0.07735"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "1.847.64"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "'4.997"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "'0.6845"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'83.03"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "'0.000009969"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").Value = "'6.178"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "29.409.86"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'229.61"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'7.568"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.965"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'157.06"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1404"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.392"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.66"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.464"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.05709"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.249"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.127"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.035"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.852"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.157"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7164"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.594"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.253.39"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.788"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01806"
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9083"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.186"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.004.87"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'101.76"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.55"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.087"
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.197"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4026"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.705"
$ws.Range("E51").Value = "  -0.01%  "
